$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M2").ClearContents()
$ws.Range("H2").Value = 699.8
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 699.8
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 699.8
$ws.Range("N2").Value = -925.8

$ws.Range("H69").Value = 3385.3635
$ws.Range("I69").Value = 3373.9
$ws.Range("J69").Value = 3500
$ws.Range("K69").Value = 10121.7
$ws.Range("L69").Value = 10500
$ws.Range("M69").Value = -9247.700000000001
$ws.Range("N69").Value = -12248

$ws.Range("H72").Value = 3385.3635
$ws.Range("I72").Value = 3373.9
$ws.Range("J72").Value = 3500
$ws.Range("K72").Value = 30365.1
$ws.Range("L72").Value = 31500
$ws.Range("M72").Value = -25997.1
$ws.Range("N72").Value = -40236

$ws.Range("H113").Value = 3597.087
$ws.Range("I113").Value = 2877.9167
$ws.Range("J113").Value = 4381.636
$ws.Range("K113").Value = 2877.9167
$ws.Range("L113").Value = 4381.636
$ws.Range("M113").Value = 376.0832999999998

$ws.Range("H132").Value = 4466880.5
$ws.Range("I132").Value = 2532.12
$ws.Range("J132").Value = 41669784
$ws.Range("K132").Value = 7596.36
$ws.Range("L132").Value = 125009352
$ws.Range("M132").Value = -5066.36
$ws.Range("N132").Value = -125014412

$ws.Range("H138").Value = 7940896
$ws.Range("I138").Value = 16669761
$ws.Range("J138").Value = 5564.091
$ws.Range("K138").Value = 50009283
$ws.Range("L138").Value = 16692.273
$ws.Range("M138").Value = -50004143
$ws.Range("N138").Value = -26972.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 928.8484999999999
$ws.Range("I74").Value = 927.4838999999999
$ws.Range("J74").Value = 950
$ws.Range("K74").Value = 927.4838999999999
$ws.Range("L74").Value = 950
$ws.Range("M74").Value = -53.48389999999995
$ws.Range("N74").Value = -2698

$ws.Range("H77").Value = 928.8484999999999
$ws.Range("I77").Value = 927.4838999999999
$ws.Range("J77").Value = 950
$ws.Range("K77").Value = 4637.4195
$ws.Range("L77").Value = 4750
$ws.Range("M77").Value = -269.4195
$ws.Range("N77").Value = -13486

$ws.Range("H88").Value = 1620
$ws.Range("I88").Value = 600
$ws.Range("J88").Value = 1875
$ws.Range("K88").Value = 600
$ws.Range("L88").Value = 1875
$ws.Range("M88").Value = -194
$ws.Range("N88").Value = -2687

$ws.Range("H91").Value = 1620
$ws.Range("I91").Value = 600
$ws.Range("J91").Value = 1875
$ws.Range("K91").Value = 600
$ws.Range("L91").Value = 1875
$ws.Range("M91").Value = 804
$ws.Range("N91").Value = -4683

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N12").ClearContents()
$ws.Range("H12").Value = 600
$ws.Range("I12").Value = 600
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 600
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -432

$ws.Range("H107").Value = 2017.1578
$ws.Range("I107").Value = 1524.875
$ws.Range("J107").Value = 2375.182
$ws.Range("K107").Value = 1524.875
$ws.Range("L107").Value = 2375.182
$ws.Range("M107").Value = 395.125
$ws.Range("N107").Value = -6215.182

$ws.Range("H134").Value = 2741.9067
$ws.Range("I134").Value = 1632.1964
$ws.Range("J134").Value = 6012.6313
$ws.Range("K134").Value = 4896.5892
$ws.Range("L134").Value = 18037.8939
$ws.Range("M134").Value = -2361.5892
$ws.Range("N134").Value = -23107.8939

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4299.875
$ws.Range("I86").Value = 4643.857
$ws.Range("J86").Value = 4158.2354
$ws.Range("K86").Value = 4643.857
$ws.Range("L86").Value = 4158.2354
$ws.Range("M86").Value = -3520.857
$ws.Range("N86").Value = -6404.2354

$ws.Range("H89").Value = 4299.875
$ws.Range("I89").Value = 4643.857
$ws.Range("J89").Value = 4158.2354
$ws.Range("K89").Value = 23219.285
$ws.Range("L89").Value = 20791.177
$ws.Range("M89").Value = -17603.285
$ws.Range("N89").Value = -32023.177

$ws.Range("H94").Value = 3774.2
$ws.Range("I94").Value = 2680.4443
$ws.Range("J94").Value = 4389.4375
$ws.Range("K94").Value = 2680.4443
$ws.Range("L94").Value = 4389.4375
$ws.Range("M94").Value = -2229.4443
$ws.Range("N94").Value = -5291.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N6").ClearContents()
$ws.Range("H6").Value = 149.3077
$ws.Range("I6").Value = 149.3077
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 447.9231
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -334.9231

$ws.Range("H76").Value = 4020.7693
$ws.Range("I76").Value = 2500
$ws.Range("J76").Value = 4147.5
$ws.Range("K76").Value = 7500
$ws.Range("L76").Value = 12442.5
$ws.Range("M76").Value = -7117
$ws.Range("N76").Value = -13208.5

$ws.Range("H79").Value = 4020.7693
$ws.Range("I79").Value = 2500
$ws.Range("J79").Value = 4147.5
$ws.Range("K79").Value = 7500
$ws.Range("L79").Value = 12442.5
$ws.Range("M79").Value = -6174
$ws.Range("N79").Value = -15094.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2769.05
$ws.Range("I80").Value = 2633.5833
$ws.Range("J80").Value = 2972.25
$ws.Range("K80").Value = 2633.5833
$ws.Range("L80").Value = 2972.25
$ws.Range("M80").Value = -1635.5833
$ws.Range("N80").Value = -4968.25

$ws.Range("H83").Value = 2769.05
$ws.Range("I83").Value = 2633.5833
$ws.Range("J83").Value = 2972.25
$ws.Range("K83").Value = 13167.9165
$ws.Range("L83").Value = 14861.25
$ws.Range("M83").Value = -8175.916499999999
$ws.Range("N83").Value = -24845.25

$ws.Range("H88").Value = 38295
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 38295
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 38295
$ws.Range("N88").Value = -39197

$ws.Range("H91").Value = 38295
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 38295
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 38295
$ws.Range("N91").Value = -41415

$ws.Range("H132").Value = 3692.311
$ws.Range("I132").Value = 3631.5813
$ws.Range("J132").Value = 4998
$ws.Range("K132").Value = 10894.7439
$ws.Range("L132").Value = 14994
$ws.Range("M132").Value = -8364.743899999999
$ws.Range("N132").Value = -20054

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2552.125
$ws.Range("I82").Value = 3131.75
$ws.Range("J82").Value = 1972.5
$ws.Range("K82").Value = 3131.75
$ws.Range("L82").Value = 1972.5
$ws.Range("M82").Value = -2770.75
$ws.Range("N82").Value = -2694.5

$ws.Range("H85").Value = 2552.125
$ws.Range("I85").Value = 3131.75
$ws.Range("J85").Value = 1972.5
$ws.Range("K85").Value = 3131.75
$ws.Range("L85").Value = 1972.5
$ws.Range("M85").Value = -1883.75
$ws.Range("N85").Value = -4468.5

$ws.Range("H132").Value = 15160660
$ws.Range("I132").Value = 6591.25
$ws.Range("J132").Value = 55571510
$ws.Range("K132").Value = 19773.75
$ws.Range("L132").Value = 166714530
$ws.Range("M132").Value = -17243.75
$ws.Range("N132").Value = -166719590

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 4320.3335
$ws.Range("I11").Value = 3001
$ws.Range("J11").Value = 4980
$ws.Range("K11").Value = 3001
$ws.Range("L11").Value = 4980
$ws.Range("M11").Value = -2859
$ws.Range("N11").Value = -5264

$ws.Range("H81").Value = 1821.5714
$ws.Range("I81").Value = 462.75
$ws.Range("J81").Value = 3633.3333
$ws.Range("K81").Value = 925.5
$ws.Range("L81").Value = 7266.6666
$ws.Range("M81").Value = 135.5
$ws.Range("N81").Value = -9388.6666

$ws.Range("H84").Value = 1821.5714
$ws.Range("I84").Value = 462.75
$ws.Range("J84").Value = 3633.3333
$ws.Range("K84").Value = 4627.5
$ws.Range("L84").Value = 36333.333
$ws.Range("M84").Value = 676.5
$ws.Range("N84").Value = -46941.333

$ws.Range("H122").Value = 1732.8108
$ws.Range("I122").Value = 2090.5417
$ws.Range("J122").Value = 1072.3846
$ws.Range("K122").Value = 6271.625100000001
$ws.Range("L122").Value = 3217.1538
$ws.Range("M122").Value = -3821.625100000001
